# Add NCSU sameAs URL to the synthetic data "work1" row in the CmsWork sheet.
# This inserts a new column (shifting cells right) at column R of row 2,
# then fills in the new cell with the NCSU catalog URL, and updates the
# concept reference that ended up shifted into column U.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CmsWork")

# Insert a new cell at R2, shifting R2:W2 right to S2:X2.
$ws.Range("R2").Insert(-4161)

# Populate the newly inserted cell with the NCSU "sameAs" URL.
$ws.Range("R2").Value = "https://d.lib.ncsu.edu/collections/catalog/0002030"

# The old T2 concept reference is now in U2; update it to concept:101.
$ws.Range("U2").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:101"
